$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the full A:AY row data for every affected row BEFORE any writes,
# so source data for later writes is not clobbered by earlier writes.
$rowData = @{}
$rowData[4] = $ws.Range("A4:AY4").Value2
$rowData[5] = $ws.Range("A5:AY5").Value2
$rowData[6] = $ws.Range("A6:AY6").Value2
$rowData[14] = $ws.Range("A14:AY14").Value2
$rowData[15] = $ws.Range("A15:AY15").Value2
$rowData[16] = $ws.Range("A16:AY16").Value2
$rowData[20] = $ws.Range("A20:AY20").Value2
$rowData[21] = $ws.Range("A21:AY21").Value2
$rowData[22] = $ws.Range("A22:AY22").Value2
$rowData[40] = $ws.Range("A40:AY40").Value2
$rowData[41] = $ws.Range("A41:AY41").Value2
$rowData[42] = $ws.Range("A42:AY42").Value2
$rowData[54] = $ws.Range("A54:AY54").Value2
$rowData[55] = $ws.Range("A55:AY55").Value2
$rowData[66] = $ws.Range("A66:AY66").Value2
$rowData[67] = $ws.Range("A67:AY67").Value2
$rowData[68] = $ws.Range("A68:AY68").Value2
$rowData[69] = $ws.Range("A69:AY69").Value2
$rowData[73] = $ws.Range("A73:AY73").Value2
$rowData[74] = $ws.Range("A74:AY74").Value2
$rowData[75] = $ws.Range("A75:AY75").Value2
$rowData[76] = $ws.Range("A76:AY76").Value2
$rowData[77] = $ws.Range("A77:AY77").Value2

# Write each destination row using the captured snapshot of its source row.
$ws.Range("A4:AY4").Value2 = $rowData[5]
$ws.Range("A5:AY5").Value2 = $rowData[6]
$ws.Range("A6:AY6").Value2 = $rowData[4]
$ws.Range("A14:AY14").Value2 = $rowData[16]
$ws.Range("A15:AY15").Value2 = $rowData[14]
$ws.Range("A16:AY16").Value2 = $rowData[15]
$ws.Range("A20:AY20").Value2 = $rowData[21]
$ws.Range("A21:AY21").Value2 = $rowData[22]
$ws.Range("A22:AY22").Value2 = $rowData[20]
$ws.Range("A40:AY40").Value2 = $rowData[42]
$ws.Range("A41:AY41").Value2 = $rowData[40]
$ws.Range("A42:AY42").Value2 = $rowData[41]
$ws.Range("A54:AY54").Value2 = $rowData[55]
$ws.Range("A55:AY55").Value2 = $rowData[54]
$ws.Range("A66:AY66").Value2 = $rowData[68]
$ws.Range("A67:AY67").Value2 = $rowData[69]
$ws.Range("A68:AY68").Value2 = $rowData[67]
$ws.Range("A69:AY69").Value2 = $rowData[66]
$ws.Range("A73:AY73").Value2 = $rowData[74]
$ws.Range("A74:AY74").Value2 = $rowData[75]
$ws.Range("A75:AY75").Value2 = $rowData[76]
$ws.Range("A76:AY76").Value2 = $rowData[77]
$ws.Range("A77:AY77").Value2 = $rowData[73]
